# Update Betfair Back/Lay odds values for Jogos_do_Dia sheet
# per commit: "Atualizando o arquivo XLSX"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.04
$ws.Range("G2").Value = 1000
$ws.Range("H2").Value = 1.04
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 1.01
$ws.Range("K2").Value = 1000
$ws.Range("P2").Value = 1.24
$ws.Range("Q2").Value = 1.01
$ws.Range("F3").Value = 1.37
$ws.Range("G3").Value = 980
$ws.Range("H3").Value = 1.04
$ws.Range("I3").Value = 980
$ws.Range("J3").Value = 1.37
$ws.Range("K3").Value = 980
$ws.Range("P3").Value = 1.24
$ws.Range("Q3").Value = 1.01
$ws.Range("F4").Value = 1.04
$ws.Range("G4").Value = 980
$ws.Range("H4").Value = 1.04
$ws.Range("I4").Value = 980
$ws.Range("J4").Value = 1.01
$ws.Range("K4").Value = 980
$ws.Range("P4").Value = 1.24
$ws.Range("Q4").Value = 1.01
$ws.Range("F5").Value = 1.04
$ws.Range("G5").Value = 1000
$ws.Range("H5").Value = 1.04
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 1.02
$ws.Range("K5").Value = 980
$ws.Range("M5").Value = 1.05
$ws.Range("P5").Value = 1.24
$ws.Range("Q5").Value = 1.05
$ws.Range("T5").Value = 1.05
$ws.Range("U5").Value = 1.05
$ws.Range("AA5").Value = 1000
$ws.Range("AB5").Value = 1000
$ws.Range("AC5").Value = 1000
$ws.Range("AG5").Value = 1000
$ws.Range("AI5").Value = 1000
$ws.Range("AM5").Value = 1000
$ws.Range("F6").Value = 1.04
$ws.Range("G6").Value = 1000
$ws.Range("H6").Value = 1.04
$ws.Range("I6").Value = 1000
$ws.Range("J6").Value = 1.01
$ws.Range("K6").Value = 1000
$ws.Range("P6").Value = 1.24
$ws.Range("Q6").Value = 1.01
$ws.Range("F7").Value = 1.04
$ws.Range("H7").Value = 1.04
$ws.Range("I7").Value = 1000
$ws.Range("J7").Value = 1.01
$ws.Range("K7").Value = 1000
$ws.Range("P7").Value = 1.24
$ws.Range("Q7").Value = 1.01
$ws.Range("F8").Value = 1.04
$ws.Range("G8").Value = 1000
$ws.Range("H8").Value = 1.04
$ws.Range("I8").Value = 1000
$ws.Range("J8").Value = 1.01
$ws.Range("K8").Value = 950
$ws.Range("P8").Value = 1.24
$ws.Range("Q8").Value = 1.01
$ws.Range("F9").Value = 1.04
$ws.Range("G9").Value = 1000
$ws.Range("H9").Value = 1.04
$ws.Range("J9").Value = 1.01
$ws.Range("K9").Value = 980
$ws.Range("N9").Value = 1.2
$ws.Range("O9").Value = 1.47
$ws.Range("P9").Value = 1.2
$ws.Range("Q9").Value = 1.47
$ws.Range("F10").Value = 1.04
$ws.Range("G10").Value = 1000
$ws.Range("H10").Value = 1.04
$ws.Range("I10").Value = 1000
$ws.Range("J10").Value = 1.01
$ws.Range("K10").Value = 980
$ws.Range("P10").Value = 1.25
$ws.Range("Q10").Value = 1.01
$ws.Range("F12").Value = 1.38
$ws.Range("H12").Value = 3.4
$ws.Range("I12").Value = 1000
$ws.Range("K12").Value = 3.65
$ws.Range("P12").Value = 1.22
$ws.Range("Q12").Value = 1.01
$ws.Range("F14").Value = 1.66
$ws.Range("G14").Value = 1.69
$ws.Range("H14").Value = 6
$ws.Range("I14").Value = 1000
$ws.Range("J14").Value = 3.6
$ws.Range("P14").Value = 1.24
$ws.Range("Q14").Value = 1.01
